$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.712432666666666
$ws.Range("H2").Value = 29.137298
$ws.Range("I2").Value = 0.4639063029983291
$ws.Range("J2").Value = 0.463906302998329
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7341896666666666
$ws.Range("N2").Value = 2.202569
$ws.Range("O2").Value = 0.4912907638668469
$ws.Range("P2").Value = 0.4912907638668469
$ws.Range("Q2").Value = 7.130767702062444
$ws.Range("R2").Value = 64.176909318562
$ws.Range("S2").Value = 0.2279128819626941
$ws.Range("T2").Value = 0.227912881962694

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.712432666666666
$ws.Range("H3").Value = 29.137298
$ws.Range("I3").Value = 0.4639063029983291
$ws.Range("J3").Value = 0.463906302998329
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7602199999999999
$ws.Range("N3").Value = 2.28066
$ws.Range("O3").Value = 0.5087092361331531
$ws.Range("P3").Value = 0.5087092361331531
$ws.Range("Q3").Value = 7.383585561853332
$ws.Range("R3").Value = 66.45227005667999
$ws.Range("S3").Value = 0.2359934210356351
$ws.Range("T3").Value = 0.235993421035635

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.124904999999998
$ws.Range("H4").Value = 27.37471499999999
$ws.Range("I4").Value = 0.4358435305594535
$ws.Range("J4").Value = 0.4358435305594534
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7341896666666666
$ws.Range("N4").Value = 2.202569
$ws.Range("O4").Value = 0.4912907638668469
$ws.Range("P4").Value = 0.4912907638668469
$ws.Range("Q4").Value = 6.699410960314998
$ws.Range("R4").Value = 60.29469864283499
$ws.Range("S4").Value = 0.2141259010549773
$ws.Range("T4").Value = 0.2141259010549773

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.124904999999998
$ws.Range("H5").Value = 27.37471499999999
$ws.Range("I5").Value = 0.4358435305594535
$ws.Range("J5").Value = 0.4358435305594534
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7602199999999999
$ws.Range("N5").Value = 2.28066
$ws.Range("O5").Value = 0.5087092361331531
$ws.Range("P5").Value = 0.5087092361331531
$ws.Range("Q5").Value = 6.936935279099997
$ws.Range("R5").Value = 62.43241751189998
$ws.Range("S5").Value = 0.2217176295044762
$ws.Range("T5").Value = 0.2217176295044761

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.318184
$ws.Range("H6").Value = 0.9545520000000001
$ws.Range("I6").Value = 0.01519779525677573
$ws.Range("J6").Value = 0.01519779525677573
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.7341896666666666
$ws.Range("N6").Value = 2.202569
$ws.Range("O6").Value = 0.4912907638668469
$ws.Range("P6").Value = 0.4912907638668469
$ws.Range("Q6").Value = 0.2336074048986667
$ws.Range("R6").Value = 2.102466644088
$ws.Range("S6").Value = 0.007466536440793294
$ws.Range("T6").Value = 0.007466536440793292

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.318184
$ws.Range("H7").Value = 0.9545520000000001
$ws.Range("I7").Value = 0.01519779525677573
$ws.Range("J7").Value = 0.01519779525677573
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.7602199999999999
$ws.Range("N7").Value = 2.28066
$ws.Range("O7").Value = 0.5087092361331531
$ws.Range("P7").Value = 0.5087092361331531
$ws.Range("Q7").Value = 0.24188984048
$ws.Range("R7").Value = 2.17700856432
$ws.Range("S7").Value = 0.007731258815982442
$ws.Range("T7").Value = 0.007731258815982441

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.353022
$ws.Range("H8").Value = 1.059066
$ws.Range("I8").Value = 0.01686180347577968
$ws.Range("J8").Value = 0.01686180347577968
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7341896666666666
$ws.Range("N8").Value = 2.202569
$ws.Range("O8").Value = 0.4912907638668469
$ws.Range("P8").Value = 0.4912907638668469
$ws.Range("Q8").Value = 0.259185104506
$ws.Range("R8").Value = 2.332665940554
$ws.Range("S8").Value = 0.008284048309788456
$ws.Range("T8").Value = 0.008284048309788454

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.353022
$ws.Range("H9").Value = 1.059066
$ws.Range("I9").Value = 0.01686180347577968
$ws.Range("J9").Value = 0.01686180347577968
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7602199999999999
$ws.Range("N9").Value = 2.28066
$ws.Range("O9").Value = 0.5087092361331531
$ws.Range("P9").Value = 0.5087092361331531
$ws.Range("Q9").Value = 0.26837438484
$ws.Range("R9").Value = 2.41536946356
$ws.Range("S9").Value = 0.008577755165991229
$ws.Range("T9").Value = 0.008577755165991227

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.427651
$ws.Range("H10").Value = 4.282953
$ws.Range("I10").Value = 0.06819056770966213
$ws.Range("J10").Value = 0.06819056770966211
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.7341896666666666
$ws.Range("N10").Value = 2.202569
$ws.Range("O10").Value = 0.4912907638668469
$ws.Range("P10").Value = 0.4912907638668469
$ws.Range("Q10").Value = 1.048166611806333
$ws.Range("R10").Value = 9.433499506257
$ws.Range("S10").Value = 0.03350139609859386
$ws.Range("T10").Value = 0.03350139609859384

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.427651
$ws.Range("H11").Value = 4.282953
$ws.Range("I11").Value = 0.06819056770966213
$ws.Range("J11").Value = 0.06819056770966211
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.7602199999999999
$ws.Range("N11").Value = 2.28066
$ws.Range("O11").Value = 0.5087092361331531
$ws.Range("P11").Value = 0.5087092361331531
$ws.Range("Q11").Value = 1.08532884322
$ws.Range("R11").Value = 9.767959588979998
$ws.Range("S11").Value = 0.03468917161106828
$ws.Range("T11").Value = 0.03468917161106827
